# Applies:
#  1) Table style swap ({D1ECE101-...} -> {414D9346-...}) on the three
#     tables that use it (slides 14, 15, 16).
#  2) The deck-wide theme color scheme being swapped from the "Integral"
#     (Red Violet) palette to the default Office "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$oldStyleId = "{D1ECE101-092A-4B6F-B965-E7AC7043960E}"
$newStyleId = "{414D9346-A596-4990-8061-E75B01794A8C}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the presentation theme palette --------------------------------
# The deck currently uses the "Integral" (Red Violet) theme; restore the
# stock Office "Office Theme" palette on the shared theme.
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407      # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$s1 = $p.Slides.Item(1)
$cs = $s1.ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $cs.Colors($idx).RGB = $officeColors[$idx]
}
